# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" and "全部类型" worksheets, reflecting refreshed data pulled
# from bilibili at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览": rows 2-11 ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 583
$wsExhibit.Range("F3").Value = 196
$wsExhibit.Range("F4").Value = 411
$wsExhibit.Range("F5").Value = 451
$wsExhibit.Range("F6").Value = 272
$wsExhibit.Range("F7").Value = 2476
$wsExhibit.Range("F8").Value = 427
$wsExhibit.Range("F9").Value = 6527
$wsExhibit.Range("F10").Value = 175
$wsExhibit.Range("F11").Value = 421

# --- Sheet "全部类型": rows 2-6 and 9-13 (rows 7-8 belong to a
#     different category and are unaffected) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 583
$wsAll.Range("F3").Value = 196
$wsAll.Range("F4").Value = 411
$wsAll.Range("F5").Value = 451
$wsAll.Range("F6").Value = 272
$wsAll.Range("F9").Value = 2476
$wsAll.Range("F10").Value = 427
$wsAll.Range("F11").Value = 6527
$wsAll.Range("F12").Value = 175
$wsAll.Range("F13").Value = 421
